# "fixed export and fixing maps"
#
# The sheet previously carried a stray population-census caption row and
# three data columns (1989 / 2002 / 2014). This edit:
#   1. Renames the sheet from the generic "1" to "თეთრიწყარო".
#   2. Removes the obsolete "(მოსახლეობის აღწერის შედეგებით)" row.
#   3. Drops the 1989 and 2002 columns, keeping only the 2014 figure.
#   4. Leaves the cursor on A2, matching the refreshed export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Give the worksheet its proper municipality name.
$ws.Name = "თეთრიწყარო"

# 2. Delete the now-unused caption row (row 2); rows below shift up.
$ws.Rows.Item(2).Delete()

# 3. Delete the 1989 and 2002 columns (both were column B once the old
#    column B shifts left after the first delete); only the 2014 figures,
#    originally in column D, remain and move into column B.
$ws.Columns.Item(2).Delete()
$ws.Columns.Item(2).Delete()

# 4. Match the saved selection state of the refreshed workbook.
$ws.Range("A2").Select()
